# Update "Pais" sheet with refreshed COVID country data and revised timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Updated timestamp (data refresh time changed from 11:22 to 11:52)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 11:52"

# Espana (row 5)
$ws.Range("B5").Value = 135032
$ws.Range("C5").Value = 3386
$ws.Range("D5").Value = 40437
$ws.Range("E5").Value = 81540
$ws.Range("F5").Value = 6931
$ws.Range("G5").Value = 414
$ws.Range("H5").Value = 13055

# Alemania (row 7)
$ws.Range("B7").Value = 100132
$ws.Range("C7").Value = 9
$ws.Range("E7").Value = 69848

# Australia (row 24)
$ws.Range("D24").Value = 2432
$ws.Range("E24").Value = 3324
$ws.Range("F24").Value = 96

# Malasia (row 33)
$ws.Range("B33").Value = 3793
$ws.Range("C33").Value = 131
$ws.Range("D33").Value = 1241
$ws.Range("E33").Value = 2490
$ws.Range("F33").Value = 102
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 62

# Eslovenia (row 63)
$ws.Range("D63").Value = 102
$ws.Range("E63").Value = 889

# Rows 89-91: Albania is newly inserted into the ranking, pushing
# Taiwan and Afganistan down one row each with their prior data.
$ws.Range("A89").Value = "Albania"
$ws.Range("B89").Value = 377
$ws.Range("C89").Value = 16
$ws.Range("D89").Value = 116
$ws.Range("E89").Value = 240
$ws.Range("F89").Value = 7
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 21

$ws.Range("A90").Value = "Taiwan"
$ws.Range("B90").Value = 373
$ws.Range("C90").Value = 10
$ws.Range("D90").Value = 57
$ws.Range("E90").Value = 311
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 5

$ws.Range("A91").Value = "Afganistan"
$ws.Range("B91").Value = 367
$ws.Range("C91").Value = 18
$ws.Range("D91").Value = 17
$ws.Range("E91").Value = 343
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 7

# Brunei (row 121)
$ws.Range("D121").Value = 82
$ws.Range("E121").Value = 52

# Rows 142-144: Etiopia is newly inserted into the ranking, pushing
# Macao and Togo down one row each with their prior data.
$ws.Range("A142").Value = "Etiopia"
$ws.Range("B142").Value = 44
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 4
$ws.Range("E142").Value = 38
$ws.Range("F142").Value = 1
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 2

$ws.Range("A143").Value = "Macao"
$ws.Range("B143").Value = 44
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 10
$ws.Range("E143").Value = 34
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

$ws.Range("A144").Value = "Togo"
$ws.Range("B144").Value = 44
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 20
$ws.Range("E144").Value = 21
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 3
